$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategy")

$ws.Range("H2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 1
$ws.Range("O2").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("S3").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("M4").Value = 0
$ws.Range("N4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("E5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 1
$ws.Range("G6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 1
$ws.Range("H7").Value = 0
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("R7").Value = 1
$ws.Range("E8").Value = 1
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 1
$ws.Range("S8").Value = 0
$ws.Range("C9").Value = 0
$ws.Range("O9").Value = 1
$ws.Range("F10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 1
$ws.Range("Q10").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = 0
$ws.Range("O11").Value = 0
$ws.Range("Q11").Value = 1
$ws.Range("S12").Value = 1
$ws.Range("K13").Value = 0
$ws.Range("P13").Value = 1
